# Apply "Add data for 2022-09-19" update:
# - Rename sheet from "Through 2022-09-10" to "Through 2022-09-11"
# - Update header cell I1 text from "2022 (through 09-10)" to "2022 (through 09-11)"
# - Update I10 (October / 2022 column) from 48 to 58
# - Update I14 (Total / 2022 column) from 1185 to 1195

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2022-09-11"

$ws.Range("I1").Value = "2022 (through 09-11)"
$ws.Range("I10").Value = 58
$ws.Range("I14").Value = 1195
